$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: swap the shift-register part for a newer one, with updated pricing
$ws.Range("A15").Value = "568-3968-5-ND"
$ws.Range("B15").Value = "74HC595 8-bit shiftregister"
$ws.Range("C15").Value = 0.1716

# New row 16: additional part (aluminum SMD capacitor)
$ws.Range("A16").Value = "PCE3878DKR-ND"
$ws.Range("B16").Value = "Al CAP 10uF 16v 20% SMD"
$ws.Range("C16").Value = 0.11528
$ws.Range("C16").Style = "Currency"
$ws.Range("D16").Value = 2
$ws.Range("E16").Formula = "=C16*D16"

# Update the active cell selection as recorded by Excel when the edit was made
$ws.Range("D17").Select()
